# Added Hawkeye simulation for both configs on gcc benchmark
#
# Fills in the previously-blank Hawkeye (row 21) and OPTGen (row 22) result
# rows for the "gcc" benchmark on both the Config1 and Config2 sheets. Those
# rows used to have only A/B labels with H/I computed from empty C:G cells
# (hence the #DIV/0! shared-formula errors); we now populate the raw
# simulation counters so the existing shared formulas (IPC in H, MPKI in I,
# and the hit-rate formula in J22) recalculate to real numbers.

$wb = $excel.ActiveWorkbook

function Set-HawkeyeOptGenRows {
    param(
        [Parameter(Mandatory = $true)] $ws,
        [Parameter(Mandatory = $true)] [double] $TotalCycles21,
        [Parameter(Mandatory = $true)] [double] $Access21,
        [Parameter(Mandatory = $true)] [double] $Hit21,
        [Parameter(Mandatory = $true)] [double] $Miss21,
        [Parameter(Mandatory = $true)] [double] $TotalCycles22,
        [Parameter(Mandatory = $true)] [double] $Access22,
        [Parameter(Mandatory = $true)] [double] $Hit22
    )

    # Row 21 - Hawkeye: C/D/E/F are raw counters, G is a literal (Total Miss),
    # H/I are the pre-existing shared formulas which just needed real inputs.
    $ws.Range("C21").Value = 50000000
    $ws.Range("D21").Value = $TotalCycles21
    $ws.Range("E21").Value = $Access21
    $ws.Range("F21").Value = $Hit21
    $ws.Range("G21").Value = $Miss21

    # Row 22 - OPTGen: same raw counters, but G22 and J22 are live formulas.
    $ws.Range("C22").Value = 50000000
    $ws.Range("D22").Value = $TotalCycles22
    $ws.Range("E22").Value = $Access22
    $ws.Range("F22").Value = $Hit22
    $ws.Range("G22").Formula = "=E22-F22"
    $ws.Range("J22").Formula = "=F22/E22"
}

$wsConfig1 = $wb.Worksheets.Item("Config1")
Set-HawkeyeOptGenRows -ws $wsConfig1 `
    -TotalCycles21 181610645 -Access21 249544 -Hit21 20244 -Miss21 229300 `
    -TotalCycles22 181610645 -Access22 29056  -Hit22 4828

$wsConfig2 = $wb.Worksheets.Item("Config2")
Set-HawkeyeOptGenRows -ws $wsConfig2 `
    -TotalCycles21 177984613 -Access21 388409 -Hit21 20202 -Miss21 368207 `
    -TotalCycles22 177984613 -Access22 16566  -Hit22 2906

# Move the visible selection from C21 to C23 on both sheets, matching where
# the author's cursor ended up after entering the new rows.
$wsConfig1.Range("C23").Select()
$wsConfig2.Range("C23").Select()

# Config2 is the tab that was active/visible when the workbook was saved;
# re-activate it last (and nudge the viewport) so it stays the one shown.
$wsConfig2.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsConfig2.Range("C23").Select()
